$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly price-list rows were reshuffled: each row below re-applies the
# correct Fecha (col D/4), Volumen (col M/13), Precio minimo (col N/14),
# Precio maximo (col O/15), Precio promedio ponderado (col P/16) and
# Precio $/Kg (col S/19) for that row, leaving every other column as-is.
# Layout: Row, Fecha(serial), Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$flat = @(
    2, 44417, 80, 1200, 1200, 1200, 1200,
    3, 45092, 120, 2600, 2600, 2600, 2600,
    4, 45097, 90, 2600, 2600, 2600, 2600,
    5, 44438, 60, 1200, 1200, 1200, 1200,
    6, 44343, 60, 1300, 1300, 1300, 1300,
    7, 44830, 50, 2500, 2500, 2500, 2500,
    8, 44424, 50, 1200, 1200, 1200, 1200,
    9, 45113, 90, 2600, 2600, 2600, 2600,
    10, 45093, 90, 2600, 2600, 2600, 2600,
    11, 45106, 120, 2600, 2600, 2600, 2600,
    12, 44435, 130, 1300, 1300, 1300, 1300,
    13, 45079, 30, 2600, 2600, 2600, 2600,
    14, 44753, 160, 2300, 2300, 2300, 2300,
    15, 44749, 120, 2300, 2300, 2300, 2300,
    16, 45044, 150, 3500, 3500, 3500, 3500,
    17, 44418, 40, 1200, 1200, 1200, 1200,
    18, 45055, 25, 2800, 2800, 2800, 2800,
    19, 45111, 50, 2600, 2600, 2600, 2600,
    20, 45042, 25, 3500, 3500, 3500, 3500,
    21, 45104, 50, 2600, 2600, 2600, 2600,
    22, 44432, 30, 1300, 1300, 1300, 1300,
    23, 45054, 25, 2500, 2500, 2500, 2500,
    24, 45041, 80, 3500, 3500, 3500, 3500,
    25, 45090, 50, 2600, 2600, 2600, 2600,
    26, 44763, 50, 2300, 2300, 2300, 2300,
    27, 44812, 50, 2500, 2500, 2500, 2500,
    28, 44405, 50, 1200, 1200, 1200, 1200,
    29, 45075, 240, 3200, 3200, 3200, 3200,
    30, 44811, 60, 2500, 2500, 2500, 2500,
    31, 44476, 80, 1200, 1200, 1200, 1200,
    32, 45086, 30, 2600, 2600, 2600, 2600,
    33, 45112, 50, 2600, 2600, 2600, 2600,
    34, 44473, 120, 1200, 1200, 1200, 1200,
    35, 44748, 300, 2300, 2300, 2300, 2300,
    36, 44357, 35, 1000, 1000, 1000, 1000,
    37, 45076, 100, 2600, 2600, 2600, 2600,
    38, 44762, 50, 2300, 2300, 2300, 2300,
    39, 45148, 280, 2750, 2750, 2750, 2750,
    40, 45062, 60, 3200, 3200, 3200, 3200,
    41, 45068, 50, 3250, 3250, 3250, 3250,
    42, 44760, 80, 2300, 2300, 2300, 2300,
    43, 44431, 100, 1300, 1300, 1300, 1300,
    44, 45149, 100, 2700, 2700, 2700, 2700,
    45, 45099, 200, 2600, 2600, 2600, 2600,
    46, 45085, 40, 2600, 2600, 2600, 2600
)

for ($i = 0; $i -lt $flat.Count; $i += 7) {
    $r = $flat[$i]
    $ws.Cells.Item($r, 4).Value  = $flat[$i + 1]   # D - Fecha
    $ws.Cells.Item($r, 13).Value = $flat[$i + 2]   # M - Volumen
    $ws.Cells.Item($r, 14).Value = $flat[$i + 3]   # N - Precio minimo
    $ws.Cells.Item($r, 15).Value = $flat[$i + 4]   # O - Precio maximo
    $ws.Cells.Item($r, 16).Value = $flat[$i + 5]   # P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value = $flat[$i + 6]   # S - Precio $/Kg
}

Write-Output "Updated $($flat.Count / 7) rows"
